$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.146235
$ws.Range("H2").Value = 0.438705
$ws.Range("I2").Value = 0.0224838618501081
$ws.Range("J2").Value = 0.0224838618501081
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1795736666666667
$ws.Range("N2").Value = 0.538721
$ws.Range("O2").Value = 0.06721938218475064
$ws.Range("P2").Value = 0.06721938218475064
$ws.Range("Q2").Value = 0.026259955145
$ws.Range("R2").Value = 0.236339596305
$ws.Range("S2").Value = 0.001511351302691551
$ws.Range("T2").Value = 0.001511351302691551
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.146235
$ws.Range("H3").Value = 0.438705
$ws.Range("I3").Value = 0.0224838618501081
$ws.Range("J3").Value = 0.0224838618501081
$ws.Range("O3").Value = 0.3682385515018647
$ws.Range("P3").Value = 0.3682385515018648
$ws.Range("Q3").Value = 0.143856244
$ws.Range("R3").Value = 1.294706196
$ws.Range("S3").Value = 0.008279424719851844
$ws.Range("T3").Value = 0.008279424719851845
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.146235
$ws.Range("H4").Value = 0.438705
$ws.Range("I4").Value = 0.0224838618501081
$ws.Range("J4").Value = 0.0224838618501081
$ws.Range("M4").Value = 1.481553
$ws.Range("N4").Value = 4.444659
$ws.Range("O4").Value = 0.5545861995390778
$ws.Range("P4").Value = 0.554586199539078
$ws.Range("Q4").Value = 0.216654902955
$ws.Range("R4").Value = 1.949894126595
$ws.Range("S4").Value = 0.01246923949441311
$ws.Range("T4").Value = 0.01246923949441311
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.146235
$ws.Range("H5").Value = 0.438705
$ws.Range("I5").Value = 0.0224838618501081
$ws.Range("J5").Value = 0.0224838618501081
$ws.Range("M5").Value = 0.02659666666666667
$ws.Range("N5").Value = 0.07979
$ws.Range("O5").Value = 0.009955866774306651
$ws.Range("P5").Value = 0.009955866774306652
$ws.Range("Q5").Value = 0.00388936355
$ws.Range("R5").Value = 0.03500427195
$ws.Range("S5").Value = 0.0002238463331515921
$ws.Range("T5").Value = 0.0002238463331515922
$ws.Range("G6").Value = 3.793107666666666
$ws.Range("I6").Value = 0.5831962851568996
$ws.Range("J6").Value = 0.5831962851568997
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1795736666666667
$ws.Range("N6").Value = 0.538721
$ws.Range("O6").Value = 0.06721938218475064
$ws.Range("P6").Value = 0.06721938218475064
$ws.Range("Q6").Value = 0.6811422517647777
$ws.Range("R6").Value = 6.130280265883
$ws.Range("S6").Value = 0.03920209398068845
$ws.Range("T6").Value = 0.03920209398068845
$ws.Range("G7").Value = 3.793107666666666
$ws.Range("I7").Value = 0.5831962851568996
$ws.Range("J7").Value = 0.5831962851568997
$ws.Range("O7").Value = 0.3682385515018647
$ws.Range("P7").Value = 0.3682385515018648
$ws.Range("Q7").Value = 3.731406448622222
$ws.Range("R7").Value = 33.5826580376
$ws.Range("S7").Value = 0.2147553552874452
$ws.Range("T7").Value = 0.2147553552874452
$ws.Range("G8").Value = 3.793107666666666
$ws.Range("I8").Value = 0.5831962851568996
$ws.Range("J8").Value = 0.5831962851568997
$ws.Range("M8").Value = 1.481553
$ws.Range("N8").Value = 4.444659
$ws.Range("O8").Value = 0.5545861995390778
$ws.Range("P8").Value = 0.554586199539078
$ws.Range("Q8").Value = 5.619690042872999
$ws.Range("R8").Value = 50.57721038585699
$ws.Range("S8").Value = 0.3234326113704732
$ws.Range("T8").Value = 0.3234326113704734
$ws.Range("G9").Value = 3.793107666666666
$ws.Range("I9").Value = 0.5831962851568996
$ws.Range("J9").Value = 0.5831962851568997
$ws.Range("M9").Value = 0.02659666666666667
$ws.Range("N9").Value = 0.07979
$ws.Range("O9").Value = 0.009955866774306651
$ws.Range("P9").Value = 0.009955866774306652
$ws.Range("Q9").Value = 0.1008840202411111
$ws.Range("R9").Value = 0.90795618217
$ws.Range("S9").Value = 0.005806224518292644
$ws.Range("T9").Value = 0.005806224518292645
$ws.Range("G10").Value = 2.288493
$ws.Range("H10").Value = 6.865479000000001
$ws.Range("I10").Value = 0.3518594075080483
$ws.Range("J10").Value = 0.3518594075080483
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1795736666666667
$ws.Range("N10").Value = 0.538721
$ws.Range("O10").Value = 0.06721938218475064
$ws.Range("P10").Value = 0.06721938218475064
$ws.Range("Q10").Value = 0.410953079151
$ws.Range("R10").Value = 3.698577712359
$ws.Range("S10").Value = 0.02365177198858342
$ws.Range("T10").Value = 0.02365177198858342
$ws.Range("G11").Value = 2.288493
$ws.Range("H11").Value = 6.865479000000001
$ws.Range("I11").Value = 0.3518594075080483
$ws.Range("J11").Value = 0.3518594075080483
$ws.Range("O11").Value = 0.3682385515018647
$ws.Range("P11").Value = 0.3682385515018648
$ws.Range("Q11").Value = 2.2512668472
$ws.Range("R11").Value = 20.2614016248
$ws.Range("S11").Value = 0.1295681985530681
$ws.Range("T11").Value = 0.1295681985530681
$ws.Range("G12").Value = 2.288493
$ws.Range("H12").Value = 6.865479000000001
$ws.Range("I12").Value = 0.3518594075080483
$ws.Range("J12").Value = 0.3518594075080483
$ws.Range("M12").Value = 1.481553
$ws.Range("N12").Value = 4.444659
$ws.Range("O12").Value = 0.5545861995390778
$ws.Range("P12").Value = 0.554586199539078
$ws.Range("Q12").Value = 3.390523669629
$ws.Range("R12").Value = 30.514713026661
$ws.Range("S12").Value = 0.1951363715819602
$ws.Range("T12").Value = 0.1951363715819602
$ws.Range("G13").Value = 2.288493
$ws.Range("H13").Value = 6.865479000000001
$ws.Range("I13").Value = 0.3518594075080483
$ws.Range("J13").Value = 0.3518594075080483
$ws.Range("M13").Value = 0.02659666666666667
$ws.Range("N13").Value = 0.07979
$ws.Range("O13").Value = 0.009955866774306651
$ws.Range("P13").Value = 0.009955866774306652
$ws.Range("Q13").Value = 0.06086628549000001
$ws.Range("R13").Value = 0.54779656941
$ws.Range("S13").Value = 0.003503065384436602
$ws.Range("T13").Value = 0.003503065384436602
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2761626666666666
$ws.Range("H14").Value = 0.8284879999999999
$ws.Range("I14").Value = 0.04246044548494399
$ws.Range("J14").Value = 0.042460445484944
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1795736666666667
$ws.Range("N14").Value = 0.538721
$ws.Range("O14").Value = 0.06721938218475064
$ws.Range("P14").Value = 0.06721938218475064
$ws.Range("Q14").Value = 0.04959154264977777
$ws.Range("R14").Value = 0.4463238838479999
$ws.Range("S14").Value = 0.00285416491278722
$ws.Range("T14").Value = 0.002854164912787221
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2761626666666666
$ws.Range("H15").Value = 0.8284879999999999
$ws.Range("I15").Value = 0.04246044548494399
$ws.Range("J15").Value = 0.042460445484944
$ws.Range("O15").Value = 0.3682385515018647
$ws.Range("P15").Value = 0.3682385515018648
$ws.Range("Q15").Value = 0.2716704206222222
$ws.Range("R15").Value = 2.4450337856
$ws.Range("S15").Value = 0.01563557294149967
$ws.Range("T15").Value = 0.01563557294149967
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2761626666666666
$ws.Range("H16").Value = 0.8284879999999999
$ws.Range("I16").Value = 0.04246044548494399
$ws.Range("J16").Value = 0.042460445484944
$ws.Range("M16").Value = 1.481553
$ws.Range("N16").Value = 4.444659
$ws.Range("O16").Value = 0.5545861995390778
$ws.Range("P16").Value = 0.554586199539078
$ws.Range("Q16").Value = 0.4091496272879999
$ws.Range("R16").Value = 3.682346645591999
$ws.Range("S16").Value = 0.02354797709223129
$ws.Range("T16").Value = 0.0235479770922313
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.2761626666666666
$ws.Range("H17").Value = 0.8284879999999999
$ws.Range("I17").Value = 0.04246044548494399
$ws.Range("J17").Value = 0.042460445484944
$ws.Range("M17").Value = 0.02659666666666667
$ws.Range("N17").Value = 0.07979
$ws.Range("O17").Value = 0.009955866774306651
$ws.Range("P17").Value = 0.009955866774306652
$ws.Range("Q17").Value = 0.00734500639111111
$ws.Range("R17").Value = 0.06610505752
$ws.Range("S17").Value = 0.0004227305384258127
$ws.Range("T17").Value = 0.0004227305384258129
